# FIN514_HW2.docx - "Adjust dividend and coupon - prob3"
#
# The Binomial-spreadsheet paragraph states the calculated note price.
# Update the value from 936.0724 to 965.4206 (the trailing period that
# ends the sentence is part of the same find/replace so the two runs
# collapse into the single run the authored edit produced).

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "936.0724.",  # FindText
    $true,        # MatchCase
    $true,        # MatchWholeWord
    $false,       # MatchWildcards
    $false,       # MatchSoundsLike
    $false,       # MatchAllWordForms
    $true,        # Forward
    1,            # Wrap (wdFindContinue)
    $false,       # Format
    "965.4206.",  # ReplaceWith
    2             # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find '936.0724.' to replace in the document."
}
